# "app movil sin barra buscar"
# Adds new time-tracking rows across sheets (Jorge Luis, Fabio, Jhon Frey,
# Deuda) and updates the active sheet/selection state left behind in the
# workbook.

$wb = $excel.ActiveWorkbook

$wsJorge = $wb.Worksheets.Item("Jorge Luis")
$wsFabio = $wb.Worksheets.Item("Fabio")
$wsJhon  = $wb.Worksheets.Item("Jhon Frey")
$wsDeuda = $wb.Worksheets.Item("Deuda")

# --- New shared strings must be introduced in this exact order so the
# --- resulting sharedStrings.xml table lines up with the authored file:
# ---   40 = "ninguno"
# ---   41 = "indicar en el boton el total en pedidos"
# ---   42 = "crear interfaces y clases para carrito"

# Jhon Frey (sheet3): fill in row 12, which previously only carried styles.
# (row 12 cells already carry date/fill styles from the template, so only
# the values need to be written.)
$wsJhon.Range("A12").Value = 42072
$wsJhon.Range("B12").Value = "Localizacion en mapas android"
$wsJhon.Range("C12").Value = 6
$wsJhon.Range("D12").Value = "3 horas"
$wsJhon.Range("E12").Value = "ninguno"
$wsJhon.Range("F12").Value = 40
$wsJhon.PageSetup.Orientation = 1

# Jorge Luis (sheet1): four new task rows.
$wsJorge.Range("A14").Value = 42072
$wsJorge.Range("A14").NumberFormat = "m/d/yy"
$wsJorge.Range("B14").Value = "indicar en el boton el total en pedidos"
$wsJorge.Range("C14").Value = 3
$wsJorge.Range("F14").Value = 0

$wsJorge.Range("A15").Value = 42072
$wsJorge.Range("A15").NumberFormat = "m/d/yy"
$wsJorge.Range("B15").Value = "crear interfaces y clases para carrito"
$wsJorge.Range("C15").Value = 8
$wsJorge.Range("F15").Value = 0

$wsJorge.Range("A16").Value = 42072
$wsJorge.Range("A16").NumberFormat = "m/d/yy"
$wsJorge.Range("B16").Value = "indicar en el boton el total en pedidos"
$wsJorge.Range("C16").Value = 3
$wsJorge.Range("F16").Value = 100

$wsJorge.Range("A17").Value = 42072
$wsJorge.Range("A17").NumberFormat = "m/d/yy"
$wsJorge.Range("B17").Value = "crear interfaces y clases para carrito"
$wsJorge.Range("C17").Value = 8
$wsJorge.Range("F17").Value = 20

# Fabio (sheet2): close out row 12 and add row 13.
$wsFabio.Range("F12").Value = 0

$wsFabio.Range("A13").Value = 41708
$wsFabio.Range("A13").NumberFormat = "m/d/yy"
$wsFabio.Range("B13").Value = "terminar portafolio"
$wsFabio.Range("C13").Value = 12
$wsFabio.Range("F13").Value = 75

# Deuda (sheet5): add row 4.
$wsDeuda.Range("A4").Value = 41707
$wsDeuda.Range("A4").NumberFormat = "m/d/yy"
$wsDeuda.Range("D4").Value = "no llego"
$wsDeuda.Range("E4").Value = 6000

# --- Selections / active sheet left by the editing session ---
$wsJhon.Activate()
$wsJhon.Range("A13").Select() | Out-Null

$wsDeuda.Activate()
$wsDeuda.Range("E6").Select() | Out-Null

$wsJorge.Activate()
$wsJorge.Range("F17").Select() | Out-Null

$wsFabio.Activate()
$wsFabio.Range("C12").Select() | Out-Null
